# This script updates the multiplication equations in the document
# to match the new set of generated values, as described by the diff.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "63×86=5418"; New = "17×67=1139" },
    @{ Old = "76×68=5168"; New = "41×55=2255" },
    @{ Old = "17×28=476";  New = "29×56=1624" },
    @{ Old = "41×25=1025"; New = "58×54=3132" },
    @{ Old = "68×29=1972"; New = "96×37=3552" },
    @{ Old = "50×18=900";  New = "34×34=1156" },
    @{ Old = "68×28=1904"; New = "50×39=1950" },
    @{ Old = "81×16=1296"; New = "75×88=6600" },
    @{ Old = "23×97=2231"; New = "18×70=1260" },
    @{ Old = "46×62=2852"; New = "47×78=3666" },
    @{ Old = "68×96=6528"; New = "82×96=7872" },
    @{ Old = "82×98=8036"; New = "57×47=2679" },
    @{ Old = "26×31=806";  New = "58×26=1508" },
    @{ Old = "26×81=2106"; New = "87×70=6090" },
    @{ Old = "23×11=253";  New = "36×58=2088" },
    @{ Old = "59×62=3658"; New = "19×62=1178" },
    @{ Old = "63×96=6048"; New = "70×99=6930" },
    @{ Old = "71×48=3408"; New = "45×35=1575" },
    @{ Old = "99×71=7029"; New = "99×59=5841" },
    @{ Old = "15×43=645";  New = "79×25=1975" },
    @{ Old = "81×54=4374"; New = "23×59=1357" },
    @{ Old = "48×62=2976"; New = "12×91=1092" },
    @{ Old = "80×13=1040"; New = "58×75=4350" },
    @{ Old = "47×36=1692"; New = "86×77=6622" },
    @{ Old = "31×65=2015"; New = "97×92=8924" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}
